# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.383.42'
$ws.Range("D3").Value = '''1.866.16'
$ws.Range("E3").Value = '''  -0.63%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '''  -0.07%  '
$ws.Range("D5").Value = '''243.31'
$ws.Range("E5").Value = '''  +0.18%  '
$ws.Range("D6").Value = '''0.7008'
$ws.Range("E6").Value = '''  -2.67%  '
$ws.Range("E7").Value = '''  -0.04%  '
$ws.Range("D8").Value = '''0.07884'
$ws.Range("E8").Value = '''  -1.66%  '
$ws.Range("D9").Value = '''0.3125'
$ws.Range("E9").Value = '''  -0.35%  '
$ws.Range("D10").Value = '''24.36'
$ws.Range("E10").Value = '''  -1.92%  '
$ws.Range("D11").Value = '''0.07783'
$ws.Range("E11").Value = '''  -4.55%  '
$ws.Range("D12").Value = '''1.873.58'
$ws.Range("E12").Value = '''  +0.00%  '
$ws.Range("D13").Value = '''5.145'
$ws.Range("E13").Value = '''  -1.44%  '
$ws.Range("D14").Value = '''92.44'
$ws.Range("E14").Value = '''  -2.38%  '
$ws.Range("D15").Value = '''0.6993'
$ws.Range("E15").Value = '''  -1.50%  '
$ws.Range("D16").Value = '''6.522'
$ws.Range("E16").Value = '''  +1.99%  '
$ws.Range("D17").Value = '''0.000008612'
$ws.Range("E17").Value = '''  +1.67%  '
$ws.Range("D18").Value = '''29.386.40'
$ws.Range("E18").Value = '''  +0.15%  '
$ws.Range("D19").Value = '''248.60'
$ws.Range("E19").Value = '''  +0.50%  '
$ws.Range("D20").Value = '''2.123.40'
$ws.Range("E20").Value = '''  +0.44%  '
$ws.Range("D21").Value = '''13.02'
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '''  -0.15%  '
$ws.Range("D23").Value = '''7.580'
$ws.Range("E23").Value = '''  -2.00%  '
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '''  -0.28%  '
$ws.Range("D25").Value = '''0.1536'
$ws.Range("E25").Value = '''  -4.29%  '
$ws.Range("D26").Value = '''8.971'
$ws.Range("E26").Value = '''  -0.74%  '
$ws.Range("D27").Value = '''160.58'
$ws.Range("E27").Value = '''  -1.26%  '
$ws.Range("D28").Value = '''18.69'
$ws.Range("E28").Value = '''  -0.96%  '
$ws.Range("D29").Value = '''1.583'
$ws.Range("E29").Value = '''  +5.23%  '
$ws.Range("D30").Value = '''4.296'
$ws.Range("E30").Value = '''  -2.45%  '
$ws.Range("D31").Value = '''4.246'
$ws.Range("E31").Value = '''  -0.81%  '
$ws.Range("D32").Value = '''1.207'
$ws.Range("E32").Value = '''  -0.67%  '
$ws.Range("D33").Value = '''0.05250'
$ws.Range("E33").Value = '''  -1.85%  '
$ws.Range("D34").Value = '''1.888'
$ws.Range("E34").Value = '''  -2.36%  '
$ws.Range("D35").Value = '''0.7586'
$ws.Range("E35").Value = '''  -0.26%  '
$ws.Range("D36").Value = '''1.183'
$ws.Range("E36").Value = '''  +0.55%  '
$ws.Range("D37").Value = '''2.706'
$ws.Range("E37").Value = '''  +0.21%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '''1.276.12'
$ws.Range("E38").Value = '''  +0.68%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01863'
$ws.Range("E39").Value = '''  -0.20%  '
$ws.Range("D40").Value = '''2.750'
$ws.Range("E40").Value = '''  -0.43%  '
$ws.Range("D41").Value = '''0.8962'
$ws.Range("E41").Value = '''  -1.22%  '
$ws.Range("D42").Value = '''109.75'
$ws.Range("E42").Value = '''  -3.16%  '
$ws.Range("D43").Value = '''5.949'
$ws.Range("E43").Value = '''  -7.62%  '
$ws.Range("D44").Value = '''70.28'
$ws.Range("E44").Value = '''  -5.12%  '
$ws.Range("D45").Value = '''1.001'
$ws.Range("E45").Value = '''  -0.04%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '''2.021.15'
$ws.Range("E46").Value = '''  -0.08%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.00000000123'
$ws.Range("E47").Value = '''  -5.01%  '
$ws.Range("D48").Value = '''9.575'
$ws.Range("E48").Value = '''  +1.01%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.791'
$ws.Range("E49").Value = '''  -0.40%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.5179'
$ws.Range("E50").Value = '''  -0.33%  '
$ws.Range("D51").Value = '''0.4285'
$ws.Range("E51").Value = '''  -1.21%  '
